$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (Artn/Gfra3 ligand-receptor pair recalculated with new TPM)
$ws.Range("I2").Value = 0.5924430993326582
$ws.Range("J2").Value = 0.5924430993326582
$ws.Range("M2").Value = 0.4802803333333334
$ws.Range("N2").Value = 1.440841
$ws.Range("Q2").Value = 0.353906570625
$ws.Range("R2").Value = 3.185159135625
$ws.Range("S2").Value = 0.5924430993326582
$ws.Range("T2").Value = 0.5924430993326582

# Update row 3 values
$ws.Range("G3").Value = 0.3936963333333334
$ws.Range("H3").Value = 1.181089
$ws.Range("I3").Value = 0.3165295008188679
$ws.Range("J3").Value = 0.3165295008188679
$ws.Range("M3").Value = 0.4802803333333334
$ws.Range("N3").Value = 1.440841
$ws.Range("Q3").Value = 0.1890846062054445
$ws.Range("R3").Value = 1.701761455849
$ws.Range("S3").Value = 0.3165295008188679
$ws.Range("T3").Value = 0.3165295008188679

# Add new row 4: Resolving-Mac -> Artn/Gfra3 -> MuSCs
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Gfra3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.113219
$ws.Range("H4").Value = 0.339657
$ws.Range("I4").Value = 0.09102739984847392
$ws.Range("J4").Value = 0.09102739984847392
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4802803333333334
$ws.Range("N4").Value = 1.440841
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.05437685905966667
$ws.Range("R4").Value = 0.489391731537
$ws.Range("S4").Value = 0.09102739984847392
$ws.Range("T4").Value = 0.09102739984847392
